$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Formula = "=C3-D3*20-E3*0.8-F3*0.6-H3*5+(U3-22)*2.5+P3/300+S3*2"
$ws.Range("N4:N8").Formula = "=C4-D4*20-E4*0.8-F4*0.6-H4*5+(U4-22)*2.5+P4/300+S4*2"
$ws.Range("N11").Formula = "=C11-D11*20-E11*0.8-F11*0.6-H11*5+(U11-22)*2.5+P11/300+S11*2"
$ws.Range("N12:N16").Formula = "=C12-D12*20-E12*0.8-F12*0.6-H12*5+(U12-22)*2.5+P12/300+S12*2"

$ws.Range("O7").Select()
